$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before J, shifting old J:L (J,K,L) to K:M.
$ws.Columns("J").EntireColumn.Insert()

# New label/input pair in row 4 (I4 label "에러코드", J4 blank input),
# matching the formatting of the other label/input pairs on that row.
$ws.Range("G4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = "에러코드"

$ws.Range("H4").Copy()
$ws.Range("J4").PasteSpecial(-4122)

# Row 5 label A5 becomes the new "서비스 ID" text; I5 is a new blank filler cell
$ws.Range("A5").Value = "서비스 ID"
$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)

# New header for inserted column J in the table header row (row 7)
$ws.Range("J7").Value = "서비스 ID"

# Update selection to match the authored state
$ws.Range("I4:J4").Select()
